# admin_console_roles.xlsx: "Teaching Fellow" role renamed to "TA".
#
# The shared-string table entry "Teaching Fellow" is replaced by a new
# entry "TA" everywhere it is used (Sheet1!C5 and backup!C5). Re-pointing
# those two cells away from the old string makes "Teaching Fellow" unused,
# so on save the writer drops it from the shared-string table and the
# other strings that used to sit after it (Designer / Teaching Staff /
# Shopper / Observer) shift down by one slot automatically - no need to
# touch C6/C10/C11/C12 directly, their displayed text is unchanged.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("backup")

$ws1.Range("C5").Value = "TA"
$ws2.Range("C5").Value = "TA"

# Selection / active-cell bookkeeping that the commit also touched.
$ws1.Activate() | Out-Null
$ws1.Range("A13").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("A1:D12").Select() | Out-Null

# Leave Sheet1 as the tab that is active/selected, matching the source.
$ws1.Activate() | Out-Null
